$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntn4"
$ws.Range("C2").Value = "Dcc"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8370283333333334
$ws.Range("H2").Value = 2.511085
$ws.Range("I2").Value = 0.01051838439934535
$ws.Range("J2").Value = 0.01051838439934535
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.092904
$ws.Range("N2").Value = 0.278712
$ws.Range("O2").Value = 0.9572435868814848
$ws.Range("P2").Value = 0.9572435868814848
$ws.Range("Q2").Value = 0.07776328028
$ws.Range("R2").Value = 0.69986952252
$ws.Range("S2").Value = 0.01006865601062759
$ws.Range("T2").Value = 0.01006865601062759

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntn4"
$ws.Range("C3").Value = "Dcc"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8370283333333334
$ws.Range("H3").Value = 2.511085
$ws.Range("I3").Value = 0.01051838439934535
$ws.Range("J3").Value = 0.01051838439934535
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.004149666666666667
$ws.Range("N3").Value = 0.012449
$ws.Range("O3").Value = 0.04275641311851519
$ws.Range("P3").Value = 0.04275641311851518
$ws.Range("Q3").Value = 0.003473388573888889
$ws.Range("R3").Value = 0.031260497165
$ws.Range("S3").Value = 0.0004497283887177551
$ws.Range("T3").Value = 0.000449728388717755

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntn4"
$ws.Range("C4").Value = "Dcc"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 34.839961
$ws.Range("H4").Value = 104.519883
$ws.Range("I4").Value = 0.437810869312907
$ws.Range("J4").Value = 0.4378108693129071
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.092904
$ws.Range("N4").Value = 0.278712
$ws.Range("O4").Value = 0.9572435868814848
$ws.Range("P4").Value = 0.9572435868814848
$ws.Range("Q4").Value = 3.236771736744
$ws.Range("R4").Value = 29.130945630696
$ws.Range("S4").Value = 0.4190916469167881
$ws.Range("T4").Value = 0.4190916469167882

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntn4"
$ws.Range("C5").Value = "Dcc"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 34.839961
$ws.Range("H5").Value = 104.519883
$ws.Range("I5").Value = 0.437810869312907
$ws.Range("J5").Value = 0.4378108693129071
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.004149666666666667
$ws.Range("N5").Value = 0.012449
$ws.Range("O5").Value = 0.04275641311851519
$ws.Range("P5").Value = 0.04275641311851518
$ws.Range("Q5").Value = 0.1445742248296667
$ws.Range("R5").Value = 1.301168023467
$ws.Range("S5").Value = 0.01871922239611892
$ws.Range("T5").Value = 0.01871922239611892

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Ntn4"
$ws.Range("C6").Value = "Dcc"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 43.90066266666667
$ws.Range("H6").Value = 131.701988
$ws.Range("I6").Value = 0.5516707462877476
$ws.Range("J6").Value = 0.5516707462877476
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.092904
$ws.Range("N6").Value = 0.278712
$ws.Range("O6").Value = 0.9572435868814848
$ws.Range("P6").Value = 0.9572435868814848
$ws.Range("Q6").Value = 4.078547164384
$ws.Range("R6").Value = 36.706924479456
$ws.Range("S6").Value = 0.5280832839540691
$ws.Range("T6").Value = 0.5280832839540691

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Ntn4"
$ws.Range("C7").Value = "Dcc"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 43.90066266666667
$ws.Range("H7").Value = 131.701988
$ws.Range("I7").Value = 0.5516707462877476
$ws.Range("J7").Value = 0.5516707462877476
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.004149666666666667
$ws.Range("N7").Value = 0.012449
$ws.Range("O7").Value = 0.04275641311851519
$ws.Range("P7").Value = 0.04275641311851518
$ws.Range("Q7").Value = 0.1821731165124445
$ws.Range("R7").Value = 1.639558048612
$ws.Range("S7").Value = 0.02358746233367852
$ws.Range("T7").Value = 0.02358746233367851
